# Re-process the sheet metadata with the newly curated dimensions.
# The "municipio" column (E) metadata block is updated:
#   E2: iaest-measure:municipio-nombre  -> sdmx-dimension:refArea
#   E3: medida                         -> dim
#   E4: xsd:int                        -> URI-Municipio

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("E3").Value = "dim"
$ws.Range("E4").Value = "URI-Municipio"
